$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 values (A7:D7) -> all become 2
$ws.Range("A7:D7").Value = 2

# Move/update the active selection from F8 to E7
$ws.Range("E7").Select()
